$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 819.6667
$ws.Range("I33").Value = 775.8
$ws.Range("K33").Value = 775.8
$ws.Range("M33").Value = -546.8

$ws.Range("H131").Value = 2257.6
$ws.Range("I131").Value = 2257.6
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6772.799999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1732.799999999999
$ws.Range("N131").ClearContents()

$ws.Range("H137").Value = 11906689
$ws.Range("J137").Value = 2048.8096
$ws.Range("L137").Value = 6146.4288
$ws.Range("N137").Value = -11246.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2221.4167
$ws.Range("I122").Value = 2221.4167
$ws.Range("K122").Value = 6664.250100000001
$ws.Range("M122").Value = -4214.250100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1933
$ws.Range("I134").Value = 1899.5
$ws.Range("K134").Value = 5698.5
$ws.Range("M134").Value = -3163.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 680.5
$ws.Range("I2").Value = 1451.75
$ws.Range("K2").Value = 1451.75
$ws.Range("M2").Value = -1338.75

$ws.Range("H3").Value = 1111
$ws.Range("I3").Value = 223
$ws.Range("K3").Value = 223
$ws.Range("M3").Value = -110

$ws.Range("H4").Value = 48
$ws.Range("I4").Value = 48
$ws.Range("K4").Value = 48
$ws.Range("M4").Value = 64

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H26").Value = 6500
$ws.Range("J26").Value = 6500
$ws.Range("L26").Value = 6500
$ws.Range("N26").Value = -7074

$ws.Range("H31").Value = 13825.909
$ws.Range("I31").Value = 14565.111
$ws.Range("K31").Value = 14565.111
$ws.Range("M31").Value = -14270.111

$ws.Range("H34").Value = 13825.909
$ws.Range("I34").Value = 14565.111
$ws.Range("K34").Value = 14565.111
$ws.Range("M34").Value = -14363.111

$ws.Range("H51").Value = 48454.453
$ws.Range("I51").Value = 44714.145
$ws.Range("K51").Value = 44714.145
$ws.Range("M51").Value = -43978.145

$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("K56").Value = 1000
$ws.Range("M56").Value = -155

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H61").Value = 48454.453
$ws.Range("I61").Value = 44714.145
$ws.Range("K61").Value = 44714.145
$ws.Range("M61").Value = -44366.145

$ws.Range("H86").Value = 13584.75
$ws.Range("I86").Value = 14147
$ws.Range("K86").Value = 14147
$ws.Range("M86").Value = -13024

$ws.Range("H89").Value = 13584.75
$ws.Range("I89").Value = 14147
$ws.Range("K89").Value = 70735
$ws.Range("M89").Value = -65119

$ws.Range("H132").Value = 22237366
$ws.Range("I132").Value = 25016536
$ws.Range("K132").Value = 75049608
$ws.Range("M132").Value = -75047078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3550172
$ws.Range("I4").Value = 1897057.8
$ws.Range("K4").Value = 5691173.4
$ws.Range("M4").Value = -5691061.4

$ws.Range("H34").Value = 43917.96
$ws.Range("J34").Value = 47691.305
$ws.Range("L34").Value = 143073.915
$ws.Range("N34").Value = -143241.915

$ws.Range("H75").Value = 1178.8
$ws.Range("J75").Value = 750
$ws.Range("L75").Value = 2250
$ws.Range("N75").Value = -4246

$ws.Range("H78").Value = 1178.8
$ws.Range("J78").Value = 750
$ws.Range("L78").Value = 6750
$ws.Range("N78").Value = -16734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 169665
$ws.Range("I7").Value = 251500
$ws.Range("K7").Value = 251500
$ws.Range("M7").Value = -251388

$ws.Range("H8").Value = 169665
$ws.Range("I8").Value = 251500
$ws.Range("K8").Value = 251500
$ws.Range("M8").Value = -251361

$ws.Range("H23").Value = 9620.799999999999
$ws.Range("I23").Value = 1007
$ws.Range("K23").Value = 1007
$ws.Range("M23").Value = -784

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 44833
$ws.Range("J3").Value = 50999.6
$ws.Range("L3").Value = 50999.6
$ws.Range("N3").Value = -51223.6

$ws.Range("H14").Value = 17371.75
$ws.Range("J14").Value = 17371.75
$ws.Range("L14").Value = 17371.75
$ws.Range("N14").Value = -17715.75

$ws.Range("H15").Value = 44833
$ws.Range("J15").Value = 50999.6
$ws.Range("L15").Value = 50999.6
$ws.Range("N15").Value = -51339.6

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H21").Value = 9998
$ws.Range("I21").Value = 9998
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 9998
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -9824
$ws.Range("N21").ClearContents()

$ws.Range("I22").Value = 4944.5
$ws.Range("K22").Value = 4944.5
$ws.Range("M22").Value = -4649.5

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("I27").Value = 4944.5
$ws.Range("K27").Value = 4944.5
$ws.Range("M27").Value = -4837.5

$ws.Range("H43").Value = 41249
$ws.Range("J43").Value = 44999
$ws.Range("L43").Value = 44999
$ws.Range("N43").Value = -45385

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 8338265
$ws.Range("J3").Value = 7397.5
$ws.Range("L3").Value = 7397.5
$ws.Range("N3").Value = -7625.5

$ws.Range("H15").Value = 12495
$ws.Range("J15").Value = 12495
$ws.Range("L15").Value = 12495
$ws.Range("N15").Value = -13071

$ws.Range("H21").Value = 39996.332
$ws.Range("J21").Value = 39995
$ws.Range("L21").Value = 39995
$ws.Range("N21").Value = -40465

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H35").Value = 39996.332
$ws.Range("J35").Value = 39995
$ws.Range("L35").Value = 39995
$ws.Range("N35").Value = -40575

$ws.Range("H39").Value = 29992.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 29992.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 29992.5
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -30818.5

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H49").Value = 60000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H54").Value = 52666.332
$ws.Range("J54").Value = 57999.625
$ws.Range("L54").Value = 57999.625
$ws.Range("N54").Value = -59039.625

$ws.Range("H124").Value = 26997.75
$ws.Range("J124").Value = 26997.75
$ws.Range("L124").Value = 26997.75
$ws.Range("N124").Value = -36817.75

$ws.Range("H125").Value = 53598.2
$ws.Range("J125").Value = 53598.2
$ws.Range("L125").Value = 53598.2
$ws.Range("N125").Value = -63438.2
